$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos price (D) and 1h volume change (E) columns per latest scrape.
# Plain decimal-looking price strings need NumberFormat forced to text first,
# otherwise Excel auto-converts them to numbers and drops the trailing zero digits
# (e.g. "303.00" -> 303, "2.430" -> 2.43) or the group separators (e.g. "1.227").

$ws.Range("D2").Value = "23.327.48"
$ws.Range("E2").Value = "  +0.13%  "

$ws.Range("D3").Value = "1.627.61"
$ws.Range("E3").Value = "  +0.91%  "

$ws.Range("E4").Value = "  +0.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.00"
$ws.Range("E6").Value = "  -0.93%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.54"
$ws.Range("E8").Value = "  -1.27%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3606"
$ws.Range("E9").Value = "  -0.27%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.227"
$ws.Range("E10").Value = "  -2.91%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.004"
$ws.Range("E11").Value = "  +0.58%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08044"
$ws.Range("E12").Value = "  -1.40%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.54"
$ws.Range("E13").Value = "  -1.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.546"
$ws.Range("E14").Value = "  -1.15%  "

$ws.Range("E15").Value = "  -0.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.198"
$ws.Range("E16").Value = "  -2.47%  "

$ws.Range("D17").Value = "1.631.59"
$ws.Range("E17").Value = "  +1.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.38"
$ws.Range("E18").Value = "  -0.88%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06928"
$ws.Range("E19").Value = "  +0.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.87"
$ws.Range("E20").Value = "  -1.98%  "

$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.441"
$ws.Range("E22").Value = "  -1.68%  "

$ws.Range("D23").Value = "23.324.50"
$ws.Range("E23").Value = "  +0.13%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.66"
$ws.Range("E24").Value = "  -2.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.211"
$ws.Range("E25").Value = "  +3.32%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.430"
$ws.Range("E26").Value = "  +0.97%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.04"
$ws.Range("E27").Value = "  -0.92%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "147.99"
$ws.Range("E28").Value = "  -1.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.301"
$ws.Range("E29").Value = "  +0.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.48"
$ws.Range("E30").Value = "  -0.72%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.296"
$ws.Range("E31").Value = "  -4.42%  "

$ws.Range("D32").Value = "1.809.73"
$ws.Range("E32").Value = "  +1.24%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.750"
$ws.Range("E33").Value = "  -1.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.81"
$ws.Range("E34").Value = "  +4.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9460"
$ws.Range("E35").Value = "  -1.26%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02815"
$ws.Range("E36").Value = "  +1.59%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2519"
$ws.Range("E37").Value = "  +0.11%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.111"
$ws.Range("E38").Value = "  -0.65%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08803"
$ws.Range("E39").Value = "  +0.23%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.07128"
$ws.Range("E40").Value = "  -3.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.363"
$ws.Range("E41").Value = "  -2.51%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7009"
$ws.Range("E42").Value = "  -1.52%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.07"
$ws.Range("E43").Value = "  +1.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.28"
$ws.Range("E44").Value = "  -2.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6432"
$ws.Range("E45").Value = "  -1.74%  "

$ws.Range("E46").Value = "  +0.85%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.307"
$ws.Range("E47").Value = "  -1.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.975"
$ws.Range("E48").Value = "  -1.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07979"
$ws.Range("E49").Value = "  -0.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.203"
$ws.Range("E50").Value = "  +0.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "126.27"
$ws.Range("E51").Value = "  -4.86%  "

Write-Output "Updated cryptos price/volume data"